$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix product_name values in column A (remove "Variasi: ..." suffix that had been
# erroneously concatenated onto the product name)
$ws.Range("A2").Value = "Oil Jug tempat minyak sally 2 Liter plastik + Tutup"
$ws.Range("A3").Value = "Rak Mini segi susun 3 Flower - Rak bumbu kosmetik serbaguna"
$ws.Range("A4").Value = "Rak plastik serbaguna susun 3 lovina warna gold"
$ws.Range("A5").Value = "Silet cukur tatra original extra tajam 10 pcs - silet serbaguna"
$ws.Range("A6").Value = "Mangko sambal bakso tutup panda star plus sendok bulat oval plastik"

# Add new row 7
$ws.Range("A7").Value = "(10bks) korek batang kayu jadul gambar grosir"
$ws.Range("B7").Value = "/(10bks)-korek-batang-kayu-jadul-gambar-grosir-i.145589728.11483476549"
$ws.Range("C7").Value = "2025-05-25 13:31"
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = "Variasi: 10bks korek pelangi"
$ws.Range("F7").Value = "menarik. kemasan nya baik. asli korek. Alhamdulillah barang nya sdh sampai dengan baik,,pesan yg pelangu datang nya yg biasa,,gpp tp dpt tambahan korek gas nya"

# Add new row 8
$ws.Range("A8").Value = "(12pcs) Solet pelet kue plastik Nice spatula murah serbaguna grosir"
$ws.Range("B8").Value = "/(12pcs)-Solet-pelet-kue-plastik-Nice-spatula-murah-serbaguna-grosir-i.145589728.3613664874"
$ws.Range("C8").Value = "2025-05-25 13:31"
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = "Variasi: 12bj solet kecil"
$ws.Range("F8").Value = "warna produk sesuai deskripsi. menggunakan bahan yg bagus dan awet. sangat awet dan tahan lama. Alhamdulillah barang nya sdh sampai dengan baik dan benar,, sesuai dengan deskripsi,, pengiriman nya sangat cepat,,"

# Add new row 9
$ws.Range("A9").Value = "erus Irus motif jagung kecil kuah sayur sendok sayutr satinless steel murah berkualitas"
$ws.Range("B9").Value = "/erus-Irus-motif-jagung-kecil-kuah-sayur-sendok-sayutr-satinless-steel-murah-berkualitas-i.145589728.6042850355"
$ws.Range("C9").Value = "2025-05-25 13:21"
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = "Variasi: Irus jagung kecil"
$ws.Range("F9").Value = "baguz. imut. bagus. Oke deh bagus ,boleh deh kpan² belanja lagi di toko ini"
